# Update "想去人数" (interested-count) column F values on sheets "展览" and "全部类型"
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 56
$ws1.Range("F3").Value = 88
$ws1.Range("F4").Value = 7232
$ws1.Range("F5").Value = 264
$ws1.Range("F6").Value = 421
$ws1.Range("F7").Value = 3714
$ws1.Range("F8").Value = 307
$ws1.Range("F9").Value = 533
$ws1.Range("F11").Value = 603
$ws1.Range("F12").Value = 94

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 56
$ws4.Range("F4").Value = 88
$ws4.Range("F6").Value = 7232
$ws4.Range("F8").Value = 264
$ws4.Range("F9").Value = 421
$ws4.Range("F10").Value = 3714
$ws4.Range("F11").Value = 307
$ws4.Range("F12").Value = 533
$ws4.Range("F14").Value = 603
$ws4.Range("F15").Value = 94
